# Automated data update: refresh the timestamp column (H) on Sheet1
# for every data row (rows 2-51) from "2025-03-07 20:19:39" to
# "2025-03-07 20:27:45".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$oldTimestamp = "2025-03-07 20:19:39"
$newTimestamp = "2025-03-07 20:27:45"

for ($row = 2; $row -le 51; $row++) {
    $cell = $ws.Cells.Item($row, 8)  # Column H
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value2 = $newTimestamp
    }
}
